# Excel COM-interop edit script for artfynd/A 13987-2022.xlsx
#
# The data rows (2-11) get cyclically rotated: the record that used to sit
# 4 rows further down now lives in this row (wrapping within the 2-11
# block), e.g. old row 6 -> new row 2, old row 7 -> new row 3, ...,
# old row 2 -> new row 8, etc. Rows 12-14 and the header row are untouched.
#
# For plain text/number cells we just overwrite .Value. The Startdatum /
# Slutdatum columns (Y, AA) hold dates written as literal text
# ("yyyy-mm-dd"); a bare string assignment would get auto-parsed into a
# date serial by Excel, so we force the cell to Text format first and
# restore the "Normal" style afterwards (so no stray number format sticks
# to the cell, matching the original un-styled cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2 ---
$ws.Range("A2").Value = 89595773
$ws.Range("B2").Value = 76909
$ws.Range("E2").Value = 6437
$ws.Range("F2").Value = "Blanksvart spiklav"
$ws.Range("G2").Value = "Calicium denigratum"
$ws.Range("H2").Value = "(Vain.) Tibell"
$ws.Range("P2").Value = "Råttenhållan, Hjd"
$ws.Range("Q2").Value = 432039.8893262132
$ws.Range("R2").Value = 6839141.943619461
$ws.Range("S2").Value = 10
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2020-09-26"
$ws.Range("Y2").Style = "Normal"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2020-09-26"
$ws.Range("AA2").Style = "Normal"
$ws.Range("AW2").Value = "Erland Lindblad"
$ws.Range("AX2").Value = "Via Erland Lindblad"
$ws.Range("AY2").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# --- row 3 ---
$ws.Range("A3").Value = 89595748
$ws.Range("B3").Value = 90841
$ws.Range("E3").Value = 2079
$ws.Range("F3").Value = "Nordtagging"
$ws.Range("G3").Value = "Odonticium romellii"
$ws.Range("H3").Value = "(S.Lundell) Parmasto"
$ws.Range("P3").Value = "Råttenhållan, Hjd"
$ws.Range("Q3").Value = 432054.1601938157
$ws.Range("R3").Value = 6839142.136896571
$ws.Range("S3").Value = 10
$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value = "2020-09-26"
$ws.Range("Y3").Style = "Normal"
$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value = "2020-09-26"
$ws.Range("AA3").Style = "Normal"
$ws.Range("AO3").Value = "Grov, gammal mossöveväxt tallåga"
$ws.Range("AW3").Value = "Erland Lindblad"
$ws.Range("AX3").Value = "Via Erland Lindblad"
$ws.Range("AY3").Value = "Kontinuitetsskogar och skogar med höga naturvärden ovan och i nära anslutning till fjällnära gränsen"

# --- row 4 ---
$ws.Range("A4").Value = 103282794
$ws.Range("B4").Value = 77259
$ws.Range("E4").Value = 228912
$ws.Range("F4").Value = "Mörk kolflarnlav"
$ws.Range("G4").Value = "Carbonicola myrmecina"
$ws.Range("H4").Value = "(Ach.) Bendiksby & Timdal"
$ws.Range("Q4").Value = 432016.9683605726
$ws.Range("R4").Value = 6839233.565669493

# --- row 5 ---
$ws.Range("A5").Value = 103282792
$ws.Range("Q5").Value = 431795.6229505471
$ws.Range("R5").Value = 6839126.828563252

# --- row 6 ---
$ws.Range("A6").Value = 103282817
$ws.Range("B6").Value = 77605
$ws.Range("E6").Value = 967
$ws.Range("F6").Value = "Varglav"
$ws.Range("G6").Value = "Letharia vulpina"
$ws.Range("H6").Value = "(L.) Hue"
$ws.Range("P6").Value = "Lillån, Hjd"
$ws.Range("Q6").Value = 431951.9551120809
$ws.Range("R6").Value = 6839240.547665785
$ws.Range("S6").Value = 25
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2022-08-30"
$ws.Range("Y6").Style = "Normal"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2022-08-30"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AW6").Value = "Andreas Öster"
$ws.Range("AX6").Value = "Andreas Öster"
$ws.Range("AY6").ClearContents()

# --- row 7 ---
$ws.Range("A7").Value = 103282804
$ws.Range("B7").Value = 78098
$ws.Range("E7").Value = 6453
$ws.Range("F7").Value = "Vedskivlav"
$ws.Range("G7").Value = "Hertelidea botryosa"
$ws.Range("H7").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("P7").Value = "Lillån, Hjd"
$ws.Range("Q7").Value = 431741.5324676937
$ws.Range("R7").Value = 6839133.123530419
$ws.Range("S7").Value = 25
$ws.Range("Y7").NumberFormat = "@"
$ws.Range("Y7").Value = "2022-08-30"
$ws.Range("Y7").Style = "Normal"
$ws.Range("AA7").NumberFormat = "@"
$ws.Range("AA7").Value = "2022-08-30"
$ws.Range("AA7").Style = "Normal"
$ws.Range("AO7").ClearContents()
$ws.Range("AW7").Value = "Andreas Öster"
$ws.Range("AX7").Value = "Andreas Öster"
$ws.Range("AY7").ClearContents()

# --- row 8 ---
$ws.Range("A8").Value = 103282810
$ws.Range("B8").Value = 78072
$ws.Range("E8").Value = 229821
$ws.Range("F8").Value = "Vedflamlav"
$ws.Range("G8").Value = "Ramboldia elabens"
$ws.Range("H8").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q8").Value = 431718.1430436089
$ws.Range("R8").Value = 6839416.587044697

# --- row 9 ---
$ws.Range("A9").Value = 103282805
$ws.Range("B9").Value = 77541
$ws.Range("E9").Value = 185
$ws.Range("F9").Value = "Violettgrå tagellav"
$ws.Range("G9").Value = "Bryoria nadvornikiana"
$ws.Range("H9").Value = "(Gyeln.) Brodo & D.Hawksw."
$ws.Range("Q9").Value = 431731.287519474
$ws.Range("R9").Value = 6839479.950734794

# --- row 10 ---
$ws.Range("A10").Value = 103282803
$ws.Range("B10").Value = 78072
$ws.Range("E10").Value = 229821
$ws.Range("F10").Value = "Vedflamlav"
$ws.Range("G10").Value = "Ramboldia elabens"
$ws.Range("H10").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q10").Value = 431688.7184720396
$ws.Range("R10").Value = 6839371.588283245

# --- row 11 ---
$ws.Range("A11").Value = 103282812
$ws.Range("B11").Value = 77605
$ws.Range("E11").Value = 967
$ws.Range("F11").Value = "Varglav"
$ws.Range("G11").Value = "Letharia vulpina"
$ws.Range("H11").Value = "(L.) Hue"
$ws.Range("Q11").Value = 431710.8905027295
$ws.Range("R11").Value = 6839410.558286428

